# Update the PriceHistory sheet with the latest price-history rows.
# New rows are inserted right under the header (row 1), pushing the
# previously-newest rows (16/10/2025, 15/10/2025) further down the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 12 fresh rows above the current row 2.
$ws.Range("A2:A13").EntireRow.Insert()

# The Price Date / Price columns look numeric to Excel's auto-detection
# (dates, decimals), so force them to plain text first -- otherwise they'd
# be stored as real date serials / floating point numbers instead of the
# literal strings used throughout the rest of the sheet.
$ws.Range("A2:C13").NumberFormat = "@"

$data = @(
    @("04/11/2025", "0.952", "SGD"),
    @("03/11/2025", "0.959", "SGD"),
    @("31/10/2025", "0.956", "SGD"),
    @("30/10/2025", "0.955", "SGD"),
    @("29/10/2025", "0.958", "SGD"),
    @("28/10/2025", "0.958", "SGD"),
    @("27/10/2025", "0.959", "SGD"),
    @("24/10/2025", "0.952", "SGD"),
    @("23/10/2025", "0.949", "SGD"),
    @("22/10/2025", "0.946", "SGD"),
    @("21/10/2025", "0.948", "SGD"),
    @("17/10/2025", "0.949", "SGD")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value() = $rec[0]
    $ws.Cells.Item($row, 2).Value() = $rec[1]
    $ws.Cells.Item($row, 3).Value() = $rec[2]
    $row = $row + 1
}
